$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: title gets a new suffix
#   "PROGETTAZIONE DEL DATABASE"
#     -> "PROGETTAZIONE DEL DATABASE - TRADUZIONE VERSO IL MODELLO RELAZIONALE"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("PROGETTAZIONE DEL DATABASE", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $rng1.InsertAfter(" - TRADUZIONE VERSO IL MODELLO RELAZIONALE")
} else {
    Write-Host "WARNING: title text not found (change 1 skipped)"
}

# ---------------------------------------------------------------------------
# Change 2: add the "email" field to the user field list
#   "cognome, " -> "cognome, email, "
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("ognome, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Text = "ognome, email, "
} else {
    Write-Host "WARNING: 'ognome, ' text not found (change 2 skipped)"
}

# ---------------------------------------------------------------------------
# Change 3: quote the table name "Sottocategoria"
#   " ha un vincolo di chiave esterna con la chiave primaria della tabella Sottocategoria."
#   -> ' ha un vincolo di chiave esterna con la chiave primaria della tabella "Sottocategoria".'
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("ha un vincolo di chiave esterna con la chiave primaria della tabella Sottocategoria.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Text = 'ha un vincolo di chiave esterna con la chiave primaria della tabella "Sottocategoria".'
} else {
    Write-Host "WARNING: foreign-key sentence not found (change 3 skipped)"
}
